$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Router" row (row 4) was edited: the quantity-per-piano (F4) was
# cleared out, and the "Numero Totale" cell (G4) - previously the
# formula =F4*8 - was overwritten with a plain literal value of 3.
$ws.Range("F4").Value = $null
$ws.Range("G4").Value = 3

# Cursor/selection ends up on F5, as in the saved file.
$ws.Range("F5").Select() | Out-Null
